$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Darshan - unchanged name, new email
$ws.Range("B2").Value = "Darshan"
$ws.Range("C2").Value = "pasne.d@husky.neu.edu"

# Row 3: Saman - unchanged name, new email
$ws.Range("B3").Value = "Saman"
$ws.Range("C3").Value = "sood.s@husky.neu.edu"

# Row 4: Shail - unchanged name, new email
$ws.Range("B4").Value = "Shail"
$ws.Range("C4").Value = "shail@ccs.neu.edu"

# Row 5: Vaibhav - unchanged name, new email
$ws.Range("B5").Value = "Vaibhav"
$ws.Range("C5").Value = "dave.v@husky.neu.edu"

# Row 6: John - unchanged name, new email
$ws.Range("B6").Value = "John"
$ws.Range("C6").Value = "snow.j@husky.neu.edu"

# Rows 7-10: new students replacing placeholder d/e/f/g entries.
# Names set first (for all four rows), then emails (for all four rows).
$ws.Range("B7").Value = "Danny"
$ws.Range("B8").Value = "Erica"
$ws.Range("B9").Value = "Flurry"
$ws.Range("B10").Value = "Gara"

$ws.Range("C7").Value = "danny.d@husky.neu.edu"
$ws.Range("C8").Value = "sniper.e@husky.neu.edu"
$ws.Range("C9").Value = "majin.f@husky.neu.edu"
$ws.Range("C10").Value = "hawking.g@husky.neu.edu"

# Rows 11-12: unchanged values (already present in shared strings)
$ws.Range("B11").Value = "Max"
$ws.Range("C11").Value = "max@x.com"

$ws.Range("B12").Value = "Kat"
$ws.Range("C12").Value = "kat@x.com"

# Add individual mailto hyperlinks for the cells whose email text no longer
# matches the shared group hyperlink's display text (C5 keeps using the
# existing C4:C12 group hyperlink).
$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:shail@ccs.neu.edu")
$ws.Hyperlinks.Add($ws.Range("C6"), "mailto:snow.j@husky.neu.edu")
$ws.Hyperlinks.Add($ws.Range("C7"), "mailto:danny.d@husky.neu.edu")
$ws.Hyperlinks.Add($ws.Range("C8"), "mailto:sniper.e@husky.neu.edu")
$ws.Hyperlinks.Add($ws.Range("C9"), "mailto:majin.f@husky.neu.edu")
$ws.Hyperlinks.Add($ws.Range("C10"), "mailto:hawking.g@husky.neu.edu")

# Keep these cells on the same "Hyperlink" cell style used by the rest of
# column C (Hyperlinks.Add re-applies formatting to the linked cell).
$ws.Range("C4").Style = "Hyperlink"
$ws.Range("C6").Style = "Hyperlink"
$ws.Range("C7").Style = "Hyperlink"
$ws.Range("C8").Style = "Hyperlink"
$ws.Range("C9").Style = "Hyperlink"
$ws.Range("C10").Style = "Hyperlink"

# Update selection to match the target
$ws.Range("C16").Select()
